# Update the German "incorrect" practice feedback text (row with key INCORRECT)
# Columns: A=key, B=EN, C=DE, D=DE_F, E=IT
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C14").Value = "Leider falsch."
$ws.Range("D14").Value = "Leider falsch."

# Match the saved selection state from the authored edit
$ws.Range("C15").Select()
